# Adds two new test-case blocks ("Instance ID Verification" and "Password
# Recovery") to the "Test Plan Final" sheet, following the same layout as
# the existing "Login"/"Sign up" blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Block 1: "Instance ID Verification" (rows 46-56), 1 test case (TC-001)
# ---------------------------------------------------------------------

# Title / Objective / Classification / Pre-requisite / blank spacer / header
# row, copied from the "Sign up" block (rows 29-35, which also contains the
# first data row we'll overwrite below).
$ws.Range("A29:G35").Copy($ws.Range("A46"))
# Tester Name / Date / Time / Pass-Fail footer rows.
$ws.Range("A39:G42").Copy($ws.Range("A53"))

# The single data row pasted above (row 52) carries the "Sign up" block's
# first-row styling; re-stripe it with the plain data-row styling used by
# the "Login" block's data rows instead.
$ws.Range("A11:E11").Copy()
$ws.Range("A52:E52").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A46").Value = "Instance ID Verification"
$ws.Range("C47").Value = "Test Instance ID verification"
$ws.Range("C49").Value = "The user must be logged in "
$ws.Range("B52").Value = "1. Login to first device.               2. Login to second device"
$ws.Range("D52").Value = "The first device should log out"

$ws.Rows.Item(46).RowHeight = 20.4
$ws.Rows.Item(47).RowHeight = 15
$ws.Rows.Item(50).RowHeight = 31.8
$ws.Rows.Item(51).RowHeight = 63
$ws.Rows.Item(52).RowHeight = 54
$ws.Rows.Item(53).RowHeight = 27
$ws.Rows.Item(54).RowHeight = 34.8
$ws.Rows.Item(55).RowHeight = 25.2

# ---------------------------------------------------------------------
# Block 2: "Password Recovery" (rows 59-70), 2 test cases (TC-001, TC-002)
# ---------------------------------------------------------------------

# Title / Objective / Classification / Pre-requisite / blank spacer / header.
$ws.Range("A29:G34").Copy($ws.Range("A59"))
# Two data rows (re-uses the "Sign up" block's first two data rows as a
# starting point for values/merges, restyled below).
$ws.Range("A35:G36").Copy($ws.Range("A65"))
# Tester Name / Date / Time / Pass-Fail footer rows.
$ws.Range("A39:G42").Copy($ws.Range("A67"))

# Re-stripe the two data rows with the plain data-row styling.
$ws.Range("A11:E11").Copy()
$ws.Range("A65:E66").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A59").Value = "Password Recovery"
$ws.Range("C60").Value = "Test if the user can recover password and forms work perfectly."
$ws.Range("C62").Value = "The user should have previously signed up"
$ws.Range("B65").Value = "1. Enter invalid email                      2. Click forget password"
$ws.Range("D65").Value = "There should be a error message and no verificationemail sent over to client."
$ws.Range("B66").Value = "1. click forget password                 2. Enter password1   in current password field             "
$ws.Range("D66").Value = ""

$ws.Rows.Item(59).RowHeight = 20.4
$ws.Rows.Item(60).RowHeight = 15
$ws.Rows.Item(65).RowHeight = 82.2
$ws.Rows.Item(66).RowHeight = 56.4

Write-Host "Added Instance ID verification and Password Recovery test cases"
